$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enter newly-graded scores -------------------------------------------
# Column D = "Odev 1" (assignment 1), Column F = "Odev 2" (assignment 2).
# Columns I (RESULT) / K (Total Vize) already hold live formulas, so they
# recalculate automatically once the raw scores below are entered.

$ws.Range("D3").Value = 65

$ws.Range("D13").Value = 55
$ws.Range("F13").Value = 55

$ws.Range("D21").Value = 55
$ws.Range("F21").Value = 50

$ws.Range("F38").Value = 60

$ws.Range("D48").Value = 60
$ws.Range("F48").Value = 60

$ws.Range("D65").Value = 50
$ws.Range("F65").Value = 55

$ws.Range("F85").Value = 70

$ws.Range("D87").Value = 0
$ws.Range("F87").Value = 45

$ws.Range("D89").Value = 50

$ws.Range("D109").Value = 60

# --- Update the sheet view / selection ------------------------------------
# Re-freeze just the header row (ySplit = 1) and leave the final selection
# on F28, matching where the author ended up after entering the scores.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F28").Select()
